$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country-name labels (column A) that moved position in the source list ---
$ws.Range('A48').Value = 'Israel'
$ws.Range('A49').Value = 'Irak'
$ws.Range('A108').Value = 'Mauritania'
$ws.Range('A109').Value = 'Costa Rica'
$ws.Range('A137').Value = 'Yemen'
$ws.Range('A138').Value = 'Uganda'
$ws.Range('A139').Value = 'San Marino'
$ws.Range('A140').Value = 'Santo Tome y Principe'
$ws.Range('A141').Value = 'Malta'
$ws.Range('A166').Value = 'Comoras'
$ws.Range('A167').Value = 'Gibraltar'
$ws.Range('A168').Value = 'Guadalupe'
$ws.Range('A169').Value = 'Siria'
$ws.Range('A184').Value = 'Botsuana'
$ws.Range('A185').Value = 'Polinesia Francesa'
$ws.Range('A206').Value = 'Groenlandia'
$ws.Range('A207').Value = 'Islas Malvinas'

# --- Update the "last updated" timestamp header ---
$ws.Range('A1').Value = 'Datos actualizados a 13 de Junio de 2020 a las 21:24'

# --- Update refreshed case-count data (columns B-H) ---
$ws.Range('B4').Value = 2132321
$ws.Range('C4').Value = 15399
$ws.Range('D4').Value = 846471
$ws.Range('E4').Value = 1168674
$ws.Range('G4').Value = 351
$ws.Range('H4').Value = 117176
$ws.Range('B7').Value = 321588
$ws.Range('C7').Value = 11985
$ws.Range('D7').Value = 162326
$ws.Range('E7').Value = 150057
$ws.Range('B12').Value = 187376
$ws.Range('C12').Value = 125
$ws.Range('E12').Value = 6611
$ws.Range('G12').Value = 2
$ws.Range('H12').Value = 8865
$ws.Range('B16').Value = 156813
$ws.Range('C16').Value = 526
$ws.Range('D16').Value = 72808
$ws.Range('E16').Value = 54607
$ws.Range('G16').Value = 24
$ws.Range('H16').Value = 29398
$ws.Range('B48').Value = 18972
$ws.Range('C48').Value = 177
$ws.Range('D48').Value = 15357
$ws.Range('E48').Value = 3315
$ws.Range('G48').Value = 0
$ws.Range('H48').Value = 300
$ws.Range('B49').Value = 18950
$ws.Range('C49').Value = 1180
$ws.Range('D49').Value = 7515
$ws.Range('E49').Value = 10886
$ws.Range('G49').Value = 53
$ws.Range('H49').Value = 549
$ws.Range('B77').Value = 4966
$ws.Range('C77').Value = 97
$ws.Range('E77').Value = 1073
$ws.Range('B107').Value = 1693
$ws.Range('C107').Value = 23
$ws.Range('D107').Value = 49
$ws.Range('E107').Value = 1617
$ws.Range('G107').Value = 3
$ws.Range('H107').Value = 27
$ws.Range('B108').Value = 1682
$ws.Range('C108').Value = 110
$ws.Range('D108').Value = 311
$ws.Range('E108').Value = 1288
$ws.Range('G108').Value = 2
$ws.Range('H108').Value = 83
$ws.Range('B109').Value = 1662
$ws.Range('C109').Value = 50
$ws.Range('D109').Value = 743
$ws.Range('E109').Value = 907
$ws.Range('H109').Value = 12
$ws.Range('B117').Value = 1357
$ws.Range('C117').Value = 36
$ws.Range('E117').Value = 243
$ws.Range('B137').Value = 705
$ws.Range('C137').Value = 73
$ws.Range('D137').Value = 39
$ws.Range('E137').Value = 506
$ws.Range('G137').Value = 21
$ws.Range('H137').Value = 160
$ws.Range('C138').Value = 8
$ws.Range('D138').Value = 219
$ws.Range('E138').Value = 475
$ws.Range('H138').Value = 0
$ws.Range('B139').Value = 694
$ws.Range('C139').Value = 0
$ws.Range('D139').Value = 520
$ws.Range('E139').Value = 132
$ws.Range('H139').Value = 42
$ws.Range('B140').Value = 650
$ws.Range('C140').Value = 11
$ws.Range('D140').Value = 168
$ws.Range('E140').Value = 470
$ws.Range('H140').Value = 12
$ws.Range('B141').Value = 646
$ws.Range('C141').Value = 1
$ws.Range('D141').Value = 601
$ws.Range('E141').Value = 36
$ws.Range('H141').Value = 9
$ws.Range('B150').Value = 486
$ws.Range('C150').Value = 14
$ws.Range('D150').Value = 247
$ws.Range('E150').Value = 236
$ws.Range('C166').Value = 13
$ws.Range('D166').Value = 114
$ws.Range('E166').Value = 60
$ws.Range('H166').Value = 2
$ws.Range('B167').Value = 176
$ws.Range('D167').Value = 173
$ws.Range('E167').Value = 3
$ws.Range('H167').Value = 0
$ws.Range('B168').Value = 171
$ws.Range('C168').Value = 0
$ws.Range('D168').Value = 157
$ws.Range('E168').Value = 0
$ws.Range('H168').Value = 14
$ws.Range('B169').Value = 170
$ws.Range('C169').Value = 6
$ws.Range('D169').Value = 71
$ws.Range('E169').Value = 93
$ws.Range('H169').Value = 6
$ws.Range('B171').Value = 142
$ws.Range('C171').Value = 1
$ws.Range('E171').Value = 6
$ws.Range('B173').Value = 138
$ws.Range('C173').Value = 8
$ws.Range('D173').Value = 61
$ws.Range('E173').Value = 71
$ws.Range('G173').Value = 1
$ws.Range('H173').Value = 6
$ws.Range('C184').Value = 12
$ws.Range('D184').Value = 24
$ws.Range('E184').Value = 35
$ws.Range('H184').Value = 1
$ws.Range('B185').Value = 60
$ws.Range('D185').Value = 60
$ws.Range('E185').Value = 0
$ws.Range('H185').Value = 0
